# "tidy away address code into functions"
# Sheet1 keeps only the first 4 address entries (rows 2-5); rows 6-11 become
# blank placeholder rows (still styled). A small bordered "input box" outline
# (red thin border) is drawn around the new D1:E11 columns next to the table.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

$xlPasteFormats = -4122
$xlNone = -4142
$xlContinuous = 1

# 1) Remove the last 6 address rows' text, keep their existing style (s=2).
$ws1.Range("A6:C11").ClearContents()

# 2) Build the red-outline box around D1:E11.
#    Use a donor cell (Sheet2!A1) that already carries fillId=2 / borderId=1
#    (white solid fill, thin colored border on all sides) so the fill is
#    reused exactly; then trim away the edges we do not want on each cell.

# D1 - top-left corner: keep left + top
$ws2.Range("A1").Copy()
$ws1.Range("D1").PasteSpecial($xlPasteFormats)
$ws1.Range("D1").Borders.Item(10).LineStyle = $xlNone   # right
$ws1.Range("D1").Borders.Item(9).LineStyle = $xlNone    # bottom

# E1 - top-right corner: keep top + right
$ws2.Range("A1").Copy()
$ws1.Range("E1").PasteSpecial($xlPasteFormats)
$ws1.Range("E1").Borders.Item(7).LineStyle = $xlNone    # left
$ws1.Range("E1").Borders.Item(9).LineStyle = $xlNone    # bottom

# D2 - left edge only (template for D2:D10)
$ws2.Range("A1").Copy()
$ws1.Range("D2").PasteSpecial($xlPasteFormats)
$ws1.Range("D2").Borders.Item(8).LineStyle = $xlNone    # top
$ws1.Range("D2").Borders.Item(10).LineStyle = $xlNone   # right
$ws1.Range("D2").Borders.Item(9).LineStyle = $xlNone    # bottom
$ws1.Range("D2").Copy()
$ws1.Range("D3:D10").PasteSpecial($xlPasteFormats)

# E2 - right edge only (template for E2:E10)
$ws2.Range("A1").Copy()
$ws1.Range("E2").PasteSpecial($xlPasteFormats)
$ws1.Range("E2").Borders.Item(8).LineStyle = $xlNone    # top
$ws1.Range("E2").Borders.Item(7).LineStyle = $xlNone    # left
$ws1.Range("E2").Borders.Item(9).LineStyle = $xlNone    # bottom
$ws1.Range("E2").Copy()
$ws1.Range("E3:E10").PasteSpecial($xlPasteFormats)

# D11 - bottom-left corner: keep left + bottom
$ws2.Range("A1").Copy()
$ws1.Range("D11").PasteSpecial($xlPasteFormats)
$ws1.Range("D11").Borders.Item(8).LineStyle = $xlNone   # top
$ws1.Range("D11").Borders.Item(10).LineStyle = $xlNone  # right

# E11 - bottom-right corner: keep right + bottom
$ws2.Range("A1").Copy()
$ws1.Range("E11").PasteSpecial($xlPasteFormats)
$ws1.Range("E11").Borders.Item(8).LineStyle = $xlNone   # top
$ws1.Range("E11").Borders.Item(7).LineStyle = $xlNone   # left

$ws1.Range("A1").Select()
